$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A10").Value = '${row.order_id}'
$ws.Range("B10").Value = '${row.city_name}'
$ws.Range("C10").Value = '${row.item_name}'
$ws.Range("D10").Value = '${row.order_date}'
$ws.Range("E10").Value = '${row.volume}'

$ws.Range("E11").Select()
